$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.378.69'
$ws.Range("E2").Value = '  -4.77%  '
$ws.Range("D3").Value = '3.011.58'
$ws.Range("E3").Value = '  -5.08%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''572.74'
$ws.Range("E5").Value = '  -4.22%  '
$ws.Range("D6").Value = '''126.32'
$ws.Range("E6").Value = '  -6.95%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").Value = '3.007.30'
$ws.Range("E8").Value = '  -5.21%  '
$ws.Range("E9").Value = '  -2.47%  '
$ws.Range("E10").Value = '  -7.04%  '
$ws.Range("D11").Value = '''5.08'
$ws.Range("E11").Value = '  -5.09%  '
$ws.Range("D12").Value = '''0.444'
$ws.Range("E12").Value = '  -2.65%  '
$ws.Range("D13").Value = '''0.0000223'
$ws.Range("E13").Value = '  -7.14%  '
$ws.Range("D14").Value = '''32.70'
$ws.Range("E14").Value = '  -5.74%  '
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("D16").Value = '3.508.38'
$ws.Range("E16").Value = '  -4.96%  '
$ws.Range("D17").Value = '3.012.52'
$ws.Range("E17").Value = '  -4.98%  '
$ws.Range("D18").Value = '60.363.54'
$ws.Range("E18").Value = '  -4.86%  '
$ws.Range("D19").Value = '''6.54'
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("D20").Value = '''431.70'
$ws.Range("E20").Value = '  -6.74%  '
$ws.Range("D21").Value = '''13.21'
$ws.Range("E21").Value = '  -5.38%  '
$ws.Range("E22").Value = '  -3.15%  '
$ws.Range("E23").Value = '  -7.74%  '
$ws.Range("D24").Value = '''12.95'
$ws.Range("E24").Value = '  -2.06%  '
$ws.Range("D25").Value = '''79.71'
$ws.Range("E25").Value = '  -4.16%  '
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").Value = '''2.55'
$ws.Range("E28").Value = '  -5.42%  '
$ws.Range("E29").Value = '  -4.55%  '
$ws.Range("D30").Value = '''7.26'
$ws.Range("E30").Value = '  -6.12%  '
$ws.Range("D31").Value = '''6.17'
$ws.Range("E31").Value = '  -9.47%  '
$ws.Range("E32").Value = '  -6.96%  '
$ws.Range("D33").Value = '''0.0963'
$ws.Range("E33").Value = '  -4.83%  '
$ws.Range("D34").Value = '''5.64'
$ws.Range("E34").Value = '  -4.35%  '
$ws.Range("D35").Value = '''0.937'
$ws.Range("E35").Value = '  -8.31%  '
$ws.Range("D36").Value = '''50.35'
$ws.Range("E36").Value = '  -2.23%  '
$ws.Range("D37").Value = '''2.06'
$ws.Range("E37").Value = '  -14.52%  '
$ws.Range("D38").Value = '''8.55'
$ws.Range("E38").Value = '  +5.21%  '
$ws.Range("D39").Value = '0.0₃0668'
$ws.Range("E39").Value = '  -9.74%  '
$ws.Range("D40").Value = '''0.0359'
$ws.Range("E40").Value = '  -8.05%  '
$ws.Range("D41").Value = '''0.108'
$ws.Range("E41").Value = '  -4.36%  '
$ws.Range("D42").Value = '''373.17'
$ws.Range("E42").Value = '  -4.95%  '
$ws.Range("D43").Value = '2.684.09'
$ws.Range("E43").Value = '  -3.81%  '
$ws.Range("D44").Value = '''2.48'
$ws.Range("E44").Value = '  -6.06%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").Value = '''121.75'
$ws.Range("E46").Value = '  -4.66%  '
$ws.Range("E47").Value = '  -6.47%  '
$ws.Range("E48").Value = '  -5.45%  '
$ws.Range("E49").Value = '  -3.29%  '
$ws.Range("D50").Value = '''23.53'
$ws.Range("E50").Value = '  -6.18%  '
$ws.Range("D51").Value = '''0.133'
$ws.Range("E51").Value = '  -1.69%  '
